# Apply updated cryptocurrency price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.297.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.359.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.47%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.25"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.64"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.53%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.26%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.76%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.14%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.65%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.43"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.22%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.781.94"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.316.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.40%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.370.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.70%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "330.23"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.52%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.46"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +15.04%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.38%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.52%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.35"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +11.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0748"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.59%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.28"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.47%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.76%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.03%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.47%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.920"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.48%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.04"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.49%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.78"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.15%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.68"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.388"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.79%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "290.45"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.79%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.65"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.14%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0220"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.43%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.19"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.55%  "
